# Atualizando o arquivo XLSX
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("Q8").Value = 1.53

# Row 13
$ws.Range("G13").Value = 1.67
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 4.33
$ws.Range("J13").Value = 2.3
$ws.Range("K13").Value = 2.1
$ws.Range("Q13").Value = 2.05
$ws.Range("R13").Value = 1.75
$ws.Range("U13").Value = 1.44
$ws.Range("V13").Value = 2.63
$ws.Range("W13").Value = 2
$ws.Range("X13").Value = 1.73
$ws.Range("AB13").Value = 12
$ws.Range("AC13").Value = 15
$ws.Range("AD13").Value = 29
$ws.Range("AE13").Value = 9.5
$ws.Range("AI13").Value = 11
$ws.Range("AJ13").Value = 23

# Row 16
$ws.Range("K16").Value = 1.8
$ws.Range("R16").Value = 1.36

# Row 18
$ws.Range("R18").Value = 1.67

# Row 19
$ws.Range("R19").Value = 1.75

# Row 20
$ws.Range("J20").Value = 2.63
